$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the per-block summary row from AVERAGE to MEDIAN (columns C:H)
# Block 1: rows 4:35 -> summary row 36
$ws.Range("C36").Formula = "=MEDIAN(C4:C35)"
$ws.Range("D36:H36").Formula = "=MEDIAN(D4:D35)"

# Block 2: rows 41:72 -> summary row 73
$ws.Range("C73").Formula = "=MEDIAN(C41:C72)"
$ws.Range("D73:H73").Formula = "=MEDIAN(D41:D72)"

# Block 3: rows 78:109 -> summary row 110
$ws.Range("C110").Formula = "=MEDIAN(C78:C109)"
$ws.Range("D110:H110").Formula = "=MEDIAN(D78:D109)"

# Block 4: rows 115:146 -> summary row 147
$ws.Range("C147").Formula = "=MEDIAN(C115:C146)"
$ws.Range("D147:H147").Formula = "=MEDIAN(D115:D146)"

# Block 5: rows 152:183 -> summary row 184
$ws.Range("C184").Formula = "=MEDIAN(C152:C183)"
$ws.Range("D184:H184").Formula = "=MEDIAN(D152:D183)"

# Leftover manual formatting from reviewing the sheet: underline on E184
# and on the empty cell F191 where the cursor ended up.
$ws.Range("E184").Font.Underline = $true
$ws.Range("F191").Font.Underline = $true

# Cursor / view ended up scrolled down with the selection on F191.
$excel.ActiveWindow.Zoom = 85
$ws.Range("F191").Select() | Out-Null
